$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Stage 1: the paragraph that used to contain only the (now relocated)
# "_GoBack" bookmark gets filled in with the new "playgame is very long..."
# text (found right after the "Fix any style bugs..." paragraph).
# ---------------------------------------------------------------------------
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $txt = $para.Range.Text
    if ($txt.StartsWith("Fix any style bugs")) {
        $targetPara = $d.Paragraphs.Item($i + 1)
        break
    }
}
if ($targetPara -eq $null) {
    throw "Could not locate the empty bookmark paragraph after 'Fix any style bugs...'"
}

$stage1Xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">The </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>playgame</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> is very long. However, since the code is working currently</w:t></w:r><w:r><w:t xml:space="preserve"> and it’s so convoluted</w:t></w:r><w:r><w:t>, I’ve decided to separate my own code</w:t></w:r><w:r><w:t xml:space="preserve"> with low coupling</w:t></w:r><w:r><w:t xml:space="preserve"> as much as possible from the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>playgame</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> method so that I didn’t need to deal too much with it</w:t></w:r><w:r><w:t xml:space="preserve"> with high cohesion within my classes and methods.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
[void]$targetPara.Range.InsertXML($stage1Xml)

# ---------------------------------------------------------------------------
# Stage 2: the "How about the tests?" / "One of the tests fails" /
# "Fix the test..." / "Are there any questionable..." block (followed by a
# trailing empty paragraph) gets restructured: two paragraphs are removed,
# two gain bold formatting and new bodies, the bookmark moves to the new
# "UnoCard equals..." paragraph, and three new paragraphs are appended.
# ---------------------------------------------------------------------------
$startPara = $null
$endPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $txt = $para.Range.Text
    if ($txt.StartsWith("How about the tests?")) {
        $startPara = $para
    }
    if ($txt.StartsWith("Are there any questionable design decisions")) {
        $endPara = $d.Paragraphs.Item($i + 1)
    }
}

if ($startPara -eq $null -or $endPara -eq $null) {
    throw "Could not locate the 'How about the tests?' .. 'Are there any questionable...' block"
}

$rangeToReplace = $d.Range($startPara.Range.Start, $endPara.Range.End)

$stage2Xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">One of the tests fails – why? </w:t></w:r></w:p><w:p><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">The </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>UnoCard</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> equals method was checking &amp;&amp; for the cards’ ids which were different.</w:t></w:r><w:r><w:t xml:space="preserve"> But that’s changing the code. Since </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>the .contains</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> method and the equals method are doing there job correctly. The contains should return a false. Therefore, changing the asserts to false will fix the test.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Are there any questionable design decisions or places where bugs could creep in depending on how methods are used? Is the visibility of attributes and methods appropriate?</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">The </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>playgame</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> method has many if statements with convoluted structural choices. The cohesion for methods inside the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>PlayUno</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> class is very low and the responsibility is very diffused. </w:t></w:r><w:r><w:t xml:space="preserve">The visibility of attributes and methods are mostly appropriate. However, due to the high coupling between classes, there are </w:t></w:r><w:r><w:t xml:space="preserve">many occasions where </w:t></w:r><w:r><w:t xml:space="preserve">many “dots” for calling methods </w:t></w:r><w:r><w:t>is a necessity.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Third Strategy:</w:t></w:r></w:p><w:p><w:r><w:t>Simply plays the lowest pointed card in the computer’s hand.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@
[void]$rangeToReplace.InsertXML($stage2Xml)
